# Daily attendance processing - swap "Miss Dina Nasr" and "Administrator"
# ordering within the "Recorded By" (column G) cells that list both roles.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2
    if ($val -eq "Miss Dina Nasr, Administrator, Developer") {
        $cell.Value2 = "Administrator, Miss Dina Nasr, Developer"
    } elseif ($val -eq "Miss Dina Nasr, Administrator") {
        $cell.Value2 = "Administrator, Miss Dina Nasr"
    }
}
